{"js": "// Apply the same textual edits described by the diff:\n// 1) \"preform\" -> \"perform\"\n// 2) \"Keeping this in binary we should\" -> \"Keeping this in binary, we should\"\n// 3) \"...left hand most side.\" -> \"...left hand most side which still gives us the value of zero below.\"\n// 4) \"(add counter point if we did not take two\u2019s complement)\" ->\n//    \"If we did ignore using 2\u2019s complement process and added the two integers together (2 + -2), then the sum would not equal zero.\"\n\nconst body = context.document.body;\n\n// --- Fix 1: typo \"preform\" -> \"perform\" ---\nconst preformHits = body.search(\"preform the signed integer arithmetic\", { matchCase: true });\npreformHits.load(\"text\");\nawait context.sync();\nif (preformHits.items.length > 0) {\n  preformHits.items[0].insertText(\"perform the signed integer arithmetic\", \"Replace\");\n  await context.sync();\n}\n\n// --- Fix 2: add comma after \"Keeping this in binary\" ---\nconst binaryHits = body.search(\"Keeping this in binary we should\", { matchCase: true });\nbinaryHits.load(\"text\");\nawait context.sync();\nif (binaryHits.items.length > 0) {\n  binaryHits.items[0].insertText(\"Keeping this in binary, we should\", \"Replace\");\n  await context.sync();\n}\n\n// --- Fix 3: append clause to the overflow sentence ---\nconst overflowHits = body.search(\"left hand most side.\", { matchCase: true });\noverflowHits.load(\"text\");\nawait context.sync();\nif (overflowHits.items.length > 0) {\n  overflowHits.items[0].insertText(\n    \"left hand most side which still gives us the value of zero below.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- Fix 4: replace the placeholder counter-point sentence ---\nconst counterPointHits = body.search(\"(add counter point if we did not take two\\u2019s complement)\", {\n  matchCase: true\n});\ncounterPointHits.load(\"text\");\nawait context.sync();\nif (counterPointHits.items.length > 0) {\n  counterPointHits.items[0].insertText(\n    \"If we did ignore using 2\\u2019s complement process and added the two integers together (2 + -2), then the sum would not equal zero.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # MatchCase=True, Forward=True, Wrap=wdFindContinue(1), Replace=wdReplaceAll(2)\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n$apos = [char]0x2019\n\n# 1) typo fix: \"preform\" -> \"perform\"\nReplace-Text \"preform the signed integer arithmetic\" \"perform the signed integer arithmetic\"\n\n# 2) add comma: \"Keeping this in binary we should\" -> \"Keeping this in binary, we should\"\nReplace-Text \"Keeping this in binary we should\" \"Keeping this in binary, we should\"\n\n# 3) extend overflow sentence with a new clause\nReplace-Text \"left hand most side.\" \"left hand most side which still gives us the value of zero below.\"\n\n# 4) replace the placeholder counter-point sentence\n$oldCounter = \"(add counter point if we did not take two\" + $apos + \"s complement)\"\n$newCounter = \"If we did ignore using 2\" + $apos + \"s complement process and added the two integers together (2 + -2), then the sum would not equal zero.\"\nReplace-Text $oldCounter $newCounter\n"}
